$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 (OverallRebateEfficiency): append Week_31..Week_40 columns (AD:AM) ---
$weekNames = @("Week_31","Week_32","Week_33","Week_34","Week_35","Week_36","Week_37","Week_38","Week_39","Week_40")
$weekValues = @(0.7347,0.862,0.8041,0.7691,0.7837,0.7951,0.8408,0.7391,0.7924,0.7745)

for ($i = 0; $i -lt 10; $i++) {
    $col = 30 + $i   # AD = column 30
    $headerCell = $ws1.Cells.Item(1, $col)
    $headerCell.Value = $weekNames[$i]
    $headerCell.Font.Size = 9

    $dataCell = $ws1.Cells.Item(2, $col)
    $dataCell.Value = $weekValues[$i]
    $dataCell.Font.Size = 9
}

# --- Sheet1 view: clear frozen/topLeft scroll, move the (inactive) selection to A3 ---
$null = $ws1.Activate()
$null = $ws1.Range("A3").Select()

# --- Sheet2 (PSA_LOLO): update the two summary values and restore it as active sheet ---
$ws2.Range("A2").Value = 38763
$ws2.Range("B2").Value = 13629

$null = $ws2.Activate()
